$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 14, shifting existing rows 14-34 down to 15-35.
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = 44629
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100112030
$ws.Range("G14").Value = "Poroto granado"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 110
$ws.Range("K14").Value = 25000
$ws.Range("L14").Value = 26000
$ws.Range("M14").Value = 25455
$ws.Range("N14").Value = "`$/saco 25 kilos"
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 1018
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
